$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to be treated as text so numeric-looking values
# (e.g. "0.9987", "44.77") are stored as strings, not auto-converted numbers,
# matching the inlineStr cells in the original workbook.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '29.103.53'
$ws.Range("E2").Value = '  -0.27%  '

$ws.Range("D3").Value = '1.830.56'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '241.75'
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("D6").Value = '0.6307'
$ws.Range("E6").Value = '  -4.72%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '44.77'
$ws.Range("E8").Value = '  +7.22%  '

$ws.Range("D9").Value = '0.07363'
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").Value = '0.2935'
$ws.Range("E10").Value = '  +0.27%  '

$ws.Range("D11").Value = '22.71'
$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("D12").Value = '0.07660'
$ws.Range("E12").Value = '  -0.91%  '

$ws.Range("D13").Value = '1.832.68'
$ws.Range("E13").Value = '  -0.68%  '

$ws.Range("D14").Value = '4.984'
$ws.Range("E14").Value = '  +0.08%  '

$ws.Range("D15").Value = '0.6641'
$ws.Range("E15").Value = '  -0.53%  '

$ws.Range("D16").Value = '82.34'
$ws.Range("E16").Value = '  -0.72%  '

$ws.Range("D17").Value = '0.000008811'
$ws.Range("E17").Value = '  +5.26%  '

$ws.Range("D18").Value = '6.069'
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("D19").Value = '29.093.58'
$ws.Range("E19").Value = '  -0.33%  '

$ws.Range("D20").Value = '2.079.98'
$ws.Range("E20").Value = '  -0.50%  '

$ws.Range("D21").Value = '226.47'
$ws.Range("E21").Value = '  -0.33%  '

$ws.Range("E22").Value = '  -0.11%  '

$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").Value = '7.159'
$ws.Range("E24").Value = '  +0.22%  '

$ws.Range("D25").Value = '0.9999'
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").Value = '158.84'
$ws.Range("E26").Value = '  -0.72%  '

$ws.Range("D27").Value = '8.459'
$ws.Range("E27").Value = '  -1.73%  '

$ws.Range("D28").Value = '0.1358'
$ws.Range("E28").Value = '  -1.95%  '

$ws.Range("D29").Value = '17.86'
$ws.Range("E29").Value = '  -0.74%  '

$ws.Range("D30").Value = '1.503'
$ws.Range("E30").Value = '  -0.46%  '

$ws.Range("D31").Value = '4.081'
$ws.Range("E31").Value = '  -0.65%  '

$ws.Range("D32").Value = '4.029'
$ws.Range("E32").Value = '  -0.42%  '

$ws.Range("D33").Value = '1.202'
$ws.Range("E33").Value = '  +2.16%  '

$ws.Range("D34").Value = '0.05318'
$ws.Range("E34").Value = '  +0.80%  '

$ws.Range("D35").Value = '1.843'
$ws.Range("E35").Value = '  -2.25%  '

$ws.Range("E36").Value = '  +2.23%  '

$ws.Range("D37").Value = '0.7341'
$ws.Range("E37").Value = '  -2.91%  '

$ws.Range("D38").Value = '2.654'
$ws.Range("E38").Value = '  -0.70%  '

$ws.Range("D39").Value = '1.299.93'
$ws.Range("E39").Value = '  +0.91%  '

$ws.Range("D40").Value = '0.01791'
$ws.Range("E40").Value = '  -0.11%  '

$ws.Range("D41").Value = '2.747'
$ws.Range("E41").Value = '  +0.92%  '

$ws.Range("D42").Value = '6.338'
$ws.Range("E42").Value = '  +6.19%  '

$ws.Range("D43").Value = '0.9007'
$ws.Range("E43").Value = '  -2.70%  '

$ws.Range("D44").Value = '0.9991'
$ws.Range("E44").Value = '  -0.69%  '

$ws.Range("D45").Value = '102.68'
$ws.Range("E45").Value = '  +0.19%  '

$ws.Range("D46").Value = '1.977.88'
$ws.Range("E46").Value = '  -0.56%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '64.54'
$ws.Range("E47").Value = '  +2.43%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '0.5124'
$ws.Range("E48").Value = '  -0.77%  '

$ws.Range("E49").Value = '  -0.29%  '

$ws.Range("D50").Value = '1.726'
$ws.Range("E50").Value = '  -2.55%  '

$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '0.3992'
$ws.Range("E51").Value = '  -1.17%  '

# Restore the default (Normal) style so no stray number-format style
# records get attached to the cells we only temporarily reformatted.
$rng.Style = "Normal"

Write-Host "Updated cryptos list"